# edit.ps1 - apply CV edits per commit "add one page versions"
$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# 1) "May 2014" (graduation date) -> " December " / "2014" (two runs)
$p2 = $d.Paragraphs(2)
$p2xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
  "<w:pPr><w:pStyle w:val='BodyText'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:firstLine='720'/></w:pPr>" + `
  "<w:r><w:t>Bachel</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'>or of Science, </w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'>Electrical </w:t></w:r>" + `
  "<w:r><w:t>Engineering</w:t></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/></w:r>" + `
  "<w:r><w:tab/><w:t xml:space='preserve'> December </w:t></w:r>" + `
  "<w:r><w:t>2014</w:t></w:r>" + `
  "</w:p>"
$p2.Range.InsertXML($p2xml)

Write-Host "Step 1 done"

# 2) Remove the first-line indent on the empty paragraph right after the
#    education GPA line (currently has <w:ind w:firstLine="720"/>).
$p6 = $d.Paragraphs(6)
$p6.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='BodyText'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/></w:pPr></w:p>")

Write-Host "Step 2 done"

# 3) NanoJapan/Rice dates: "June 2011 - July 2011" -> "June 2012 - July 2012"
$null = $d.Content.Find.Execute("     June 2011 " + [char]0x2013 + " July 2011", $false, $false, $false, $false, $false, `
    $true, $wdFindContinue, $false, "     June 2012 " + [char]0x2013 + " July 2012", $wdReplaceOne)

Write-Host "Step 3 done"

# 4) Delete the stray empty paragraph right after the "Instron..." bullet
#    (just before the "INDUSTRY" section heading).
$p28 = $d.Paragraphs(28)
$p28.Range.Delete()

Write-Host "Step 4 done"

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 5) SpaceX bullet "...Matlab, C++, and Bash" loses its trailing <w:br/> and the
#    "_GoBack" bookmark (the bookmark is relocated below, in step 8).
$pBash = $d.Paragraphs(35)
$bashXml = "<w:p $wns>" + `
  "<w:pPr><w:pStyle w:val='BodyText'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='17'/></w:numPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='1080'/><w:rPr><w:szCs w:val='20'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:t>Developed and qualified proprietary avionics systems to improve safety and reliability of all future Falcon 9 and Falcon Heavy flights, using Matlab, C++, and Bash</w:t></w:r>" + `
  "</w:p>"
$pBash.Range.InsertXML($bashXml)

Write-Host "Step 5 done"

# 6) "Engineering and Science Tutor" loses <w:lastRenderedPageBreak/> (moved to step 7)
$pTutor = $d.Paragraphs(36)
$tutorXml = "<w:p $wns>" + `
  "<w:pPr><w:pStyle w:val='BodyText'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='720'/><w:rPr><w:szCs w:val='20'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:b/><w:szCs w:val='20'/></w:rPr><w:t>Engineering and Science Tutor</w:t></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:t>, instaEDU.com</w:t></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:tab/></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:tab/></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:tab/></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:tab/></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:tab/><w:t xml:space='preserve'>                       May 2013 " + [char]0x2013 + " Present</w:t></w:r>" + `
  "</w:p>"
$pTutor.Range.InsertXML($tutorXml)

Write-Host "Step 6 done"

# 7) "Taught science, math, ..." gains <w:lastRenderedPageBreak/>
$pTaught = $d.Paragraphs(38)
$taughtXml = "<w:p $wns>" + `
  "<w:pPr><w:pStyle w:val='BodyText'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='17'/></w:numPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='1080'/><w:rPr><w:szCs w:val='20'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:lastRenderedPageBreak/><w:t>Taught science, math, and engineering concepts to students ranging in age from middle school to college</w:t></w:r>" + `
  "</w:p>"
$pTaught.Range.InsertXML($taughtXml)

Write-Host "Step 7 done"

# 8) "...solving equations" gains a trailing "s" run plus the relocated
#    "_GoBack" bookmark.
$pEquations = $d.Paragraphs(39)
$equationsXml = "<w:p $wns>" + `
  "<w:pPr><w:pStyle w:val='BodyText'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='17'/></w:numPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='1080'/><w:rPr><w:szCs w:val='20'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:t>Designed and developed a proof-of-concept math training resource to visually teach students about solving equations</w:t></w:r>" + `
  "<w:r><w:rPr><w:szCs w:val='20'/></w:rPr><w:t>s</w:t></w:r>" + `
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
  "</w:p>"
$pEquations.Range.InsertXML($equationsXml)

Write-Host "Step 8 done"
